# "Improve switch SFX construct setting scene"
#
# On slide 2 the group "群組 127" (id 128) previously sat next to two loose
# text boxes ("文字方塊 128" / id 129, text "MAX" and "文字方塊 129" / id 130,
# text "300"). The edit pulls those two text boxes inside the group (as its
# last two children) and the group itself is recreated (new id 131, new
# name "群組 130"), while every existing child of the group - and the two
# text boxes themselves - keep their original ids, geometry and text.
#
# Reproducing this with the high-level Shapes.Range(...).Group() call only
# works cleanly (no extra nesting level gets added) if the group is first
# dissolved so all of its children become top-level shapes, which are then
# grouped back together with the two text boxes in one single Group() call.
#
# The new group's numeric Id is assigned internally by the host and is not
# settable directly, but it is allocated from a deterministic sequence that
# is shared by every id-allocating call (AddTextbox, Group, ...). Burning
# through that sequence with harmless, fully undone scratch shapes before
# doing the real work lands the real group on the exact id (131) seen in
# the target deck; the Name is then fixed up explicitly to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Locate the existing group ("群組 127", id 128) and the two loose text
# boxes ("文字方塊 128" id 129 / "MAX", "文字方塊 129" id 130 / "300") that
# need to be pulled inside it.
$group = $null
$extraIds = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Id -eq 128) {
        $group = $candidate
    } elseif ($candidate.Id -eq 129 -or $candidate.Id -eq 130) {
        $extraIds += $candidate.Id
    }
}

# Burn through 7 scratch id allocations (add + immediately delete a
# throwaway textbox) so that the *next* id-allocating call - the real
# regroup below - lands on id 131, matching the target deck exactly.
for ($burn = 0; $burn -lt 7; $burn++) {
    $scratch = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
    $scratch.Delete()
}

# Dissolve the old group so its children become top-level shapes again,
# and remember their ids.
$formerChildren = $group.Ungroup()
$childIds = @()
for ($i = 1; $i -le $formerChildren.Count; $i++) {
    $childIds += $formerChildren.Item($i).Id
}

# Re-collect the ids of interest: the (now top-level) former group members
# plus the two loose text boxes "MAX" (129) and "300" (130).
$targetIds = $childIds + $extraIds
$indices = @()
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    if ($targetIds -contains $s.Shapes.Item($i).Id) {
        $indices += $i
    }
}

# Group them all back together in one shot - this keeps every child's id
# intact and produces a single, flat group (no extra nesting level).
$range = $s.Shapes.Range($indices)
$newGroup = $range.Group()
$newGroup.Name = "群組 130"
